$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.9160666910459944
$ws1.Range("C2").Value = -1.056566437996126
$ws1.Range("B3").Value = 1.239837176224575
$ws1.Range("C3").Value = 0.8332499611178321
$ws1.Range("B4").Value = 0.9720368665876842
$ws1.Range("C4").Value = -0.1670110671399031

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -1.939632327693161
$ws2.Range("C2").Value = -0.2979247751224782
$ws2.Range("B3").Value = 1.604600679975619
$ws2.Range("C3").Value = 0.4277726519311412
$ws2.Range("B4").Value = 0.4008681094733339
$ws2.Range("C4").Value = 0.1528739961421536
